# Update the metadata sheet for indicator 3.3.2 ("Add files via upload").
#
# - B4:  indicator title, drop the trailing colon after "3.3.2"
# - B6:  contact organization name, de-duplicate "КР КР" -> "КР" (with a
#        stray typo "jтдел" in place of "отдел", as in the source upload)
# - B10: organization website, "www.stat.kg" -> "www.stat.gov.kg"
#
# Re-assigning Font.Name (to the font the cell already has) mirrors what
# Excel does when a cell is retyped/edited in this workbook: it mints a
# fresh (but visually identical) cell style for the touched cells instead
# of reusing the old one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "3.3.2 Заболеваемость туберкулезом на 100 000 человек"
$ws.Range("B4").Font.Name = "Calibri"

$ws.Range("B6").Value = "Национальный статистический комитет КР (jтдел социальной статистики)"
$ws.Range("B6").Font.Name = "Calibri"

$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"

# Restore the selection to B4 (matches the saved view state in the upload).
$ws.Range("B4").Select()
